# Update the "K" column (column G) values for the save_data sheet.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the underlying data regeneration results in new K
# values for each existing row (rows 2-29) while leaving everything else intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 7
    5  = 2
    6  = 6
    7  = 7
    8  = 3
    9  = 6
    10 = 6
    11 = 4
    12 = 8
    13 = 6
    14 = 5
    15 = 8
    16 = 7
    17 = 6
    18 = 6
    19 = 4
    20 = 6
    21 = 6
    22 = 6
    23 = 6
    24 = 5
    25 = 2
    26 = 9
    27 = 10
    28 = 4
    29 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
